$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. "559.38", "1.00") but must be stored as
# literal text, matching the inlineStr cells in the original workbook. Force
# a temporary text number format before assignment so Excel does not coerce
# the string into a number, then restore the default "Normal" style so no
# extra style index is left attached to the cell.
$dCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D11", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D25", "D30", "D34", "D36", "D37", "D38", "D40", "D41", "D45", "D46", "D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.896.23"
$ws.Range("D3").Value = "3.063.03"
$ws.Range("D5").Value = "559.38"
$ws.Range("D6").Value = "142.78"
$ws.Range("D8").Value = "3.063.78"
$ws.Range("D9").Value = "0.516"
$ws.Range("D11").Value = "6.15"
$ws.Range("D14").Value = "35.37"
$ws.Range("D15").Value = "3.563.52"
$ws.Range("D16").Value = "63.948.68"
$ws.Range("D17").Value = "3.064.87"
$ws.Range("D19").Value = "6.79"
$ws.Range("D20").Value = "487.50"
$ws.Range("D21").Value = "14.36"
$ws.Range("D22").Value = "0.691"
$ws.Range("D23").Value = "14.70"
$ws.Range("D25").Value = "82.67"
$ws.Range("D30").Value = "1.00"
$ws.Range("D34").Value = "5.70"
$ws.Range("D36").Value = "54.85"
$ws.Range("D37").Value = "0.0412"
$ws.Range("D38").Value = "444.27"
$ws.Range("D40").Value = "3.047.47"
$ws.Range("D41").Value = "2.77"
$ws.Range("D45").Value = "28.04"
$ws.Range("D46").Value = "2.26"
$ws.Range("D50").Value = "117.99"

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Column E values (percentages padded with spaces) are kept as plain text
# automatically since Excel will not parse the padded/spaced strings as numbers.
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("E11").Value = "  -3.42%  "
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("E21").Value = "  +3.63%  "
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("E23").Value = "  +8.96%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("E38").Value = "  -5.14%  "
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("E41").Value = "  -6.60%  "
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("E44").Value = "  +6.64%  "
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("E46").Value = "  +5.23%  "
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("E51").Value = "  +3.35%  "
